# Apply updated "dSF" (column F) values as part of a data repull / mean
# recalculation pass. Only column F values change; all other columns and
# data remain untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8  = 4
    10 = 5
    20 = -3
    23 = -5
    25 = -4
    27 = 0
    30 = -4
    32 = -2
    35 = -2
    37 = 2
    40 = -2
    41 = 1
    42 = 5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
